# Natmi following Dr Hou advice
# Update the Jam3-Jam2 LR-pair results: the number of ligand- and
# receptor-expressing cells changed from 1 to 3 for every edge, which in
# turn changes all of the derived expression / specificity statistics
# that NATMI recomputes from that count (columns E, G-J, K, M-P, Q-T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => column/value pairs taken from the recomputed NATMI output.
$updates = @{
    2  = @{ E=3; G=16.944289;          H=50.832867;          I=0.2968502282295034; J=0.2968502282295034;
             K=3; M=91.47356466666666;  N=274.420694;         O=0.7914295280104694; P=0.7914295280104694;
             Q=1549.954515572189;       R=13949.5906401497;   S=0.234936036017476;  T=0.234936036017476 }
    3  = @{ E=3; G=16.944289;          H=50.832867;          I=0.2968502282295034; J=0.2968502282295034;
             K=3; M=20.94207233333333;  N=62.82621699999999;  O=0.181190866265331;  P=0.181190866265331;
             Q=354.8485258749043;       R=3193.636732874138;  S=0.05378655000396494; T=0.05378655000396493 }
    4  = @{ E=3; G=16.944289;          H=50.832867;          I=0.2968502282295034; J=0.2968502282295034;
             K=3; M=3.164539666666667;  N=9.493619000000001;  O=0.02737960572419959; P=0.02737960572419959;
             Q=53.62087466396367;       R=482.587871975673;   S=0.008127642208062467; T=0.008127642208062466 }
    5  = @{ E=3; G=7.875246;           H=23.625738;          I=0.1379679355364798; J=0.1379679355364798;
             K=3; M=91.47356466666666;  N=274.420694;         O=0.7914295280104694; P=0.7914295280104694;
             Q=720.376824246908;        R=6483.391418222171;  S=0.1091918981022151; T=0.1091918981022151 }
    6  = @{ E=3; G=7.875246;           H=23.625738;          I=0.1379679355364798; J=0.1379679355364798;
             K=3; M=20.94207233333333;  N=62.82621699999999;  O=0.181190866265331;  P=0.181190866265331;
             Q=164.923971374794;        R=1484.315742373146;  S=0.02499852975669412; T=0.02499852975669412 }
    7  = @{ E=3; G=7.875246;           H=23.625738;          I=0.1379679355364798; J=0.1379679355364798;
             K=3; M=3.164539666666667;  N=9.493619000000001;  O=0.02737960572419959; P=0.02737960572419959;
             Q=24.921528351758;         R=224.293755165822;   S=0.003777507677570603; T=0.003777507677570602 }
    8  = @{ E=3; G=32.26072766666666;  H=96.78218299999999;  I=0.5651818362340167; J=0.5651818362340167;
             K=3; M=91.47356466666666;  N=274.420694;         O=0.7914295280104694; P=0.7914295280104694;
             Q=2951.003758410555;       R=26559.03382569499;  S=0.4473015938907783; T=0.4473015938907783 }
    9  = @{ E=3; G=32.26072766666666;  H=96.78218299999999;  I=0.5651818362340167; J=0.5651818362340167;
             K=3; M=20.94207233333333;  N=62.82621699999999;  O=0.181190866265331;  P=0.181190866265331;
             Q=675.6064923213011;       R=6080.458430891709;  S=0.1024057865046719; T=0.1024057865046719 }
    10 = @{ E=3; G=32.26072766666666;  H=96.78218299999999;  I=0.5651818362340167; J=0.5651818362340167;
             K=3; M=3.164539666666667;  N=9.493619000000001;  O=0.02737960572419959; P=0.02737960572419959;
             Q=102.0903523766974;       R=918.8131713902769;  S=0.01547445583856652; T=0.01547445583856652 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
